# Applies the WesternSydney 6-9m_bi_GF_Trial cleanup edit:
#  - remove the documentation sheets (Sheet2, Sheet3)
#  - rename Sheet1 to "6-9m_Bil"
#  - fix the K1 header typo: "fixation_incongruent" -> "fixation_incongrent"
#  - replace numeric 0 placeholders in column G with the text "NA"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Fix header typo in K1
$ws1.Range("K1").Value = "fixation_incongrent"

# Replace numeric 0 placeholders in column G with text "NA" for the rows
# where no first-shift/latency data exists.
$naRows = @(6, 17, 18, 33, 48, 58, 75, 76, 83)
foreach ($r in $naRows) {
    $ws1.Cells.Item($r, 7).Value = "NA"
}

# Remove the documentation worksheets
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Rename the remaining sheet
$ws1.Name = "6-9m_Bil"
